$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (K), mirroring the formatting
# of the existing "2021" column (J) so that styles carry over exactly.
$ws.Range("J4:J14").Copy($ws.Range("K4:K14"))

$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.3
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Move the active selection to match the author's final cursor position.
$ws.Range("L7").Select()
